$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a "price" cell (column D) as literal text, matching the
# workbook's convention of storing prices as inline/shared strings even
# when they look numeric (e.g. "324.08"). Prefixing with an apostrophe
# forces Excel to store the literal text instead of re-parsing it as a
# number (which would also introduce floating point noise).
function Set-PriceText {
    param([int]$Row, [string]$Text)
    $ws.Cells.Item($Row, 4).Value = "'" + $Text
}

function Set-VolumeText {
    param([int]$Row, [string]$Text)
    $ws.Cells.Item($Row, 5).Value = $Text
}

# Row 2 - Bitcoin
Set-PriceText  2 "27.581.47"
Set-VolumeText 2 "  -0.63%  "

# Row 3 - Ethereum
Set-PriceText  3 "1.750.88"
Set-VolumeText 3 "  +0.16%  "

# Row 4 - TetherUSD (price unchanged)
Set-VolumeText 4 "  -0.04%  "

# Row 5 - BNB
Set-PriceText  5 "324.08"
Set-VolumeText 5 "  +1.24%  "

# Row 6 - USDC (price unchanged)
Set-VolumeText 6 "  -0.03%  "

# Row 7 - XRP
Set-PriceText  7 "0.4586"
Set-VolumeText 7 "  +8.97%  "

# Row 8 - Cardano
Set-PriceText  8 "0.3571"
Set-VolumeText 8 "  -2.03%  "

# Row 9 - Dogecoin
Set-PriceText  9 "0.07466"
Set-VolumeText 9 "  +0.95%  "

# Row 10 - OKB
Set-PriceText  10 "42.11"
Set-VolumeText 10 "  -1.98%  "

# Row 11 - Polygon
Set-PriceText  11 "1.090"
Set-VolumeText 11 "  +0.67%  "

# Row 12 - BinanceUSD (price unchanged)
Set-VolumeText 12 "  -0.10%  "

# Row 13 - Solana
Set-PriceText  13 "20.72"
Set-VolumeText 13 "  +0.92%  "

# Row 14 - Polkadot
Set-PriceText  14 "5.993"
Set-VolumeText 14 "  -0.83%  "

# Row 15 - Chainlink
Set-PriceText  15 "7.070"
Set-VolumeText 15 "  -2.57%  "

# Row 16 - WrappedEther
Set-PriceText  16 "1.753.36"
Set-VolumeText 16 "  -1.01%  "

# Row 17 - Litecoin
Set-PriceText  17 "92.32"
Set-VolumeText 17 "  +1.77%  "

# Row 18 - ShibaInu
Set-PriceText  18 "0.00001061"
Set-VolumeText 18 "  +1.23%  "

# Row 19 - TRON
Set-PriceText  19 "0.06424"
Set-VolumeText 19 "  +1.44%  "

# Row 20 - Dai (price unchanged)
Set-VolumeText 20 "  -0.12%  "

# Row 21 - Avalanche
Set-PriceText  21 "16.73"
Set-VolumeText 21 "  -1.41%  "

# Row 22 - Uniswap
Set-PriceText  22 "5.795"
Set-VolumeText 22 "  -2.26%  "

# Row 23 - WrappedBTC
Set-PriceText  23 "27.637.69"
Set-VolumeText 23 "  -0.52%  "

# Row 24 - Cosmos
Set-PriceText  24 "11.20"
Set-VolumeText 24 "  -0.04%  "

# Row 25 - Toncoin
Set-PriceText  25 "2.105"
Set-VolumeText 25 "  +1.37%  "

# Row 26 - Monero
Set-PriceText  26 "164.15"
Set-VolumeText 26 "  +4.50%  "

# Row 27 - EthereumClassic
Set-PriceText  27 "20.22"
Set-VolumeText 27 "  +0.72%  "

# Row 28 - WrappedliquidstakedEther2.0
Set-PriceText  28 "1.955.81"
Set-VolumeText 28 "  -0.64%  "

# Rows 29 & 30 swap places (LidoDAOToken <-> BitcoinCash), with new
# price/volume values for each coin in its new row.
$ws.Cells.Item(29, 2).Value = "BitcoinCash"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-PriceText  29 "126.13"
Set-VolumeText 29 "  +1.96%  "

$ws.Cells.Item(30, 2).Value = "LidoDAOToken"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-PriceText  30 "2.060"
Set-VolumeText 30 "  -3.33%  "

# Row 31 - ImmutableX
Set-PriceText  31 "1.054"
Set-VolumeText 31 "  -6.84%  "

# Row 32 - Stellar
Set-PriceText  32 "0.09191"
Set-VolumeText 32 "  +4.22%  "

# Row 33 - HuobiToken (price unchanged)
Set-VolumeText 33 "  +0.54%  "

# Row 34 - Filecoin
Set-PriceText  34 "5.516"
Set-VolumeText 34 "  -0.42%  "

# Row 35 - Aptos
Set-PriceText  35 "11.83"
Set-VolumeText 35 "  -3.27%  "

# Row 36 - VeChain
Set-PriceText  36 "0.02288"
Set-VolumeText 36 "  +1.00%  "

# Rows 37 & 38 swap places (Algorand <-> Hedera), with new price/volume
# values for each coin in its new row.
$ws.Cells.Item(37, 2).Value = "Hedera"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-PriceText  37 "0.06034"
Set-VolumeText 37 "  +0.96%  "

$ws.Cells.Item(38, 2).Value = "Algorand"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-PriceText  38 "0.2093"
Set-VolumeText 38 "  -0.24%  "

# Row 39 - InternetComputer(DFINITY)
Set-PriceText  39 "4.963"
Set-VolumeText 39 "  +0.39%  "

# Row 40 - TheSandbox
Set-PriceText  40 "0.6306"
Set-VolumeText 40 "  +0.45%  "

# Row 41 - TrustWalletToken
Set-PriceText  41 "1.209"
Set-VolumeText 41 "  +2.97%  "

# Row 42 - WEMIXTOKEN (price unchanged)
Set-VolumeText 42 "  -0.57%  "

# Row 43 - FraxShare
Set-PriceText  43 "7.747"
Set-VolumeText 43 "  -0.16%  "

# Row 44 - EnergySwap
Set-PriceText  44 "13.30"
Set-VolumeText 44 "  -0.46%  "

# Row 45 - Decentraland
Set-PriceText  45 "0.5899"
Set-VolumeText 45 "  +0.84%  "

# Row 46 - PancakeSwap
Set-PriceText  46 "3.714"
Set-VolumeText 46 "  +1.15%  "

# Row 47 - Quant
Set-PriceText  47 "122.60"
Set-VolumeText 47 "  +0.53%  "

# Row 48 - NEARProtocol
Set-PriceText  48 "1.936"
Set-VolumeText 48 "  -1.20%  "

# Row 49 - EOS
Set-PriceText  49 "1.136"
Set-VolumeText 49 "  -2.95%  "

# Row 50 - Cronos
Set-PriceText  50 "0.06858"
Set-VolumeText 50 "  +0.88%  "

# Row 51 - Aave
Set-PriceText  51 "71.82"
Set-VolumeText 51 "  -2.15%  "
